$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 544260125102889.7
$ws.Range("C2").Value = 544260125102889.7
$ws.Range("D2").Value = 544260125102889.7

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 2980449616755.818
$ws.Range("C3").Value = 14823470709594.77
$ws.Range("D3").Value = 1549366279204.544

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 548168968728.8655
$ws.Range("C4").Value = 2507760325856.404
$ws.Range("D4").Value = 704788674079.952

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 82334981463217.02
$ws.Range("C5").Value = 81646593419727.94
$ws.Range("D5").Value = 107426465640572.4
